$wb = $excel.ActiveWorkbook

# Rows in "展览" sheet (F column = "想去人数") that were refreshed.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 15339
$ws1.Range("F9").Value  = 15333
$ws1.Range("F11").Value = 8910
$ws1.Range("F14").Value = 78
$ws1.Range("F19").Value = 41
$ws1.Range("F20").Value = 536
$ws1.Range("F24").Value = 1102
$ws1.Range("F27").Value = 69
$ws1.Range("F32").Value = 36
$ws1.Range("F33").Value = 237
$ws1.Range("F37").Value = 5450

# Same rows (shifted) in "全部类型" sheet, which mirrors "展览" data.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 15339
$ws4.Range("F9").Value  = 15333
$ws4.Range("F11").Value = 8910
$ws4.Range("F15").Value = 78
$ws4.Range("F20").Value = 41
$ws4.Range("F21").Value = 536
$ws4.Range("F25").Value = 1102
$ws4.Range("F28").Value = 69
$ws4.Range("F35").Value = 36
$ws4.Range("F36").Value = 237
$ws4.Range("F40").Value = 5450
